$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "21/04/2023"
$ws.Range("B2").Value = "Brussels"
$ws.Range("C2").Value = "Gavi, the Vaccine Alliance"
$ws.Range("D2").Value = "Global Health architecture & vaccine donations"
